$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.365.95"
$ws.Range("E2").Value = "  +3.45%  "
$ws.Range("D3").Value = "2.128.44"
$ws.Range("E3").Value = "  +3.83%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.54"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.05%  "
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "2.439.87"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.795"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").Value = "2.129.42"
$ws.Range("E17").Value = "  +3.83%  "
$ws.Range("D18").Value = "38.205.32"
$ws.Range("E18").Value = "  +3.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("D21").Value = "0.0₃0829"
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.139"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.77%  "
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0629"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.74%  "
$ws.Range("E37").Value = "  +5.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("E40").Value = "  +8.94%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("D44").Value = "1.467.07"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("E46").Value = "  +5.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("D51").Value = "2.323.79"
$ws.Range("E51").Value = "  +3.75%  "
